$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# --- Row 105: "The Net and the Butterfly" ---
$ws.Cells.Item(105,1).Value = "The Net and the Bufferfly"
$ws.Cells.Item(105,2).Value = "Olivia Fox Cabane;Judah Pollack"

$ws.Cells.Item(104,3).Copy()
$ws.Cells.Item(105,3).PasteSpecial(-4122)
$ws.Cells.Item(105,3).Value = 44019

$ws.Cells.Item(104,4).Copy()
$ws.Cells.Item(105,4).PasteSpecial(-4122)
$ws.Cells.Item(105,4).Value = 44020

$ws.Cells.Item(105,5).Value = "eureka;break through;psychology;neuroscience;happiness;productivity;creativity"
$ws.Cells.Item(105,6).Value = "Audio"
$ws.Cells.Item(105,7).Value = "9 Hours 42 Mins"
$ws.Cells.Item(105,8).Value = 3
$ws.Cells.Item(105,9).Value = $true

# --- Row 106: "The Blind Side" ---
$ws.Cells.Item(106,1).Value = "The Blind Side"
$ws.Cells.Item(106,2).Value = "Michael Lewis"

$ws.Cells.Item(104,3).Copy()
$ws.Cells.Item(106,3).PasteSpecial(-4122)
$ws.Cells.Item(106,3).Value = 44021

$ws.Cells.Item(104,4).Copy()
$ws.Cells.Item(106,4).PasteSpecial(-4122)
$ws.Cells.Item(106,4).Value = 44022

$ws.Cells.Item(106,5).Value = "football;self improvement;biography;Michael Oherr;redemption"
$ws.Cells.Item(106,6).Value = "Audio"
$ws.Cells.Item(106,7).Value = "11 Hours 49 Mins"
$ws.Cells.Item(106,8).Value = 4
$ws.Cells.Item(106,9).Value = $true

# Move the active selection to reflect where the user ended up after
# entering the new rows (matches the post-edit view state).
[void]$ws.Range("A107").Select()
